# Apply weekly update: a new record was inserted at row 73 and the
# previously existing records for rows 73-82 shifted down to rows 74-83.
# The old row 83 record is dropped (rows 84+ are untouched).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D, J, K, L, M, P contain the values that change row to row.
# New values per row (after the shift), in order D,J,K,L,M,P
$data = @{
    73 = @(44795, 160, 12000, 12000, 12000, 923)
    74 = @(44508, 400, 13000, 15000, 13850, 1065)
    75 = @(44775, 200, 11000, 12000, 11400, 877)
    76 = @(44335, 170, 25000, 27000, 25824, 1986)
    77 = @(44627, 180, 14000, 15000, 14444, 1111)
    78 = @(44252, 130, 33000, 35000, 34077, 2621)
    79 = @(44315, 400, 25000, 26000, 25425, 1956)
    80 = @(44711, 380, 12000, 13000, 12605, 970)
    81 = @(44431, 260, 12000, 13000, 12462, 959)
    82 = @(44749, 170, 12000, 13000, 12412, 955)
    83 = @(44376, 580, 12000, 14000, 13103, 1008)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("D$row").Value = $vals[0]
    $ws.Range("J$row").Value = $vals[1]
    $ws.Range("K$row").Value = $vals[2]
    $ws.Range("L$row").Value = $vals[3]
    $ws.Range("M$row").Value = $vals[4]
    $ws.Range("P$row").Value = $vals[5]
}
